# Revert to last commit before v4
# On the "Commands" sheet, a row describing the new `pidSVbuttons(<bool>)`
# command needs to be inserted above row 113 ("pidRS(<rs>)"), pushing all
# following rows down by one. The sheet's Print Area (and the view's
# selection) need to follow the shift from C132 -> C133 / C113.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Make sure the Commands sheet is the active one (it already is in the
# source workbook, but be explicit).
$ws.Activate()

# Insert a new blank row at row 113; everything below (old rows 113-167)
# shifts down to 114-168, carrying its formatting (column B italic style,
# column C normal style) along automatically.
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row with the restored command description.
$ws.Range("B113").Value = "pidSVbuttons(<bool>)"
$ws.Range("C113").Value = "toggles the visibility of the PID SV buttons"

# The sheet's print area referenced the "Commands" label cell in column C;
# after inserting the row it now lives one row further down.
$ws.PageSetup.PrintArea = '$C$133'

# Restore the active selection to the (shifted) label cell.
[void]$ws.Range("C113").Select()
